$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet held yearly data for 2004-2020 in rows 2-18.
# The update drops the oldest six years (2004-2009, old rows 2-7) and
# appends a new row for 2021 at the end, leaving rows 2-13 for 2010-2021.

# Deleting rows 2:7 shifts 2010..2020 (old rows 8-18) up to rows 2-12,
# carrying their values/number formats/styles along automatically.
$ws.Rows("2:7").Delete()

# Row 13 is now blank. Clone the formatting of row 12 (2020, the former
# last row) for column A so the new year label keeps the same style
# (bold/centered/bordered, s="1"), then fill in the 2021 figures.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 87.2283
$ws.Range("C13").Value = 83673
$ws.Range("D13").Value = 21745.7
$ws.Range("E13").Value = 32459300
$ws.Range("F13").Value = 53.247643
$ws.Range("G13").Value = 1053655.18
